$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-06-19 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-06-20 Friday", 2) | Out-Null
$d.Content.Find.Execute("54÷7=7, 5", $true, $false, $false, $false, $false, $true, 1, $false, "82÷9=9, 1", 2) | Out-Null
$d.Content.Find.Execute("90÷4=22, 2", $true, $false, $false, $false, $false, $true, 1, $false, "26÷3=8, 2", 2) | Out-Null
$d.Content.Find.Execute("90÷3=30, 0", $true, $false, $false, $false, $false, $true, 1, $false, "35÷4=8, 3", 2) | Out-Null
$d.Content.Find.Execute("84÷3=28, 0", $true, $false, $false, $false, $false, $true, 1, $false, "32÷7=4, 4", 2) | Out-Null
$d.Content.Find.Execute("91÷9=10, 1", $true, $false, $false, $false, $false, $true, 1, $false, "54÷6=9, 0", 2) | Out-Null
$d.Content.Find.Execute("82÷5=16, 2", $true, $false, $false, $false, $false, $true, 1, $false, "62÷7=8, 6", 2) | Out-Null
$d.Content.Find.Execute("74÷2=37, 0", $true, $false, $false, $false, $false, $true, 1, $false, "58÷3=19, 1", 2) | Out-Null
$d.Content.Find.Execute("22÷9=2, 4", $true, $false, $false, $false, $false, $true, 1, $false, "68÷9=7, 5", 2) | Out-Null
$d.Content.Find.Execute("77÷6=12, 5", $true, $false, $false, $false, $false, $true, 1, $false, "95÷8=11, 7", 2) | Out-Null
$d.Content.Find.Execute("88÷7=12, 4", $true, $false, $false, $false, $false, $true, 1, $false, "73÷6=12, 1", 2) | Out-Null
$d.Content.Find.Execute("29÷9=3, 2", $true, $false, $false, $false, $false, $true, 1, $false, "91÷7=13, 0", 2) | Out-Null
$d.Content.Find.Execute("90÷2=45, 0", $true, $false, $false, $false, $false, $true, 1, $false, "12÷6=2, 0", 2) | Out-Null
$d.Content.Find.Execute("78÷3=26, 0", $true, $false, $false, $false, $false, $true, 1, $false, "64÷2=32, 0", 2) | Out-Null
$d.Content.Find.Execute("46÷7=6, 4", $true, $false, $false, $false, $false, $true, 1, $false, "45÷5=9, 0", 2) | Out-Null
$d.Content.Find.Execute("30÷5=6, 0", $true, $false, $false, $false, $false, $true, 1, $false, "41÷5=8, 1", 2) | Out-Null
$d.Content.Find.Execute("51÷7=7, 2", $true, $false, $false, $false, $false, $true, 1, $false, "89÷3=29, 2", 2) | Out-Null
$d.Content.Find.Execute("37÷7=5, 2", $true, $false, $false, $false, $false, $true, 1, $false, "53÷2=26, 1", 2) | Out-Null
$d.Content.Find.Execute("32÷9=3, 5", $true, $false, $false, $false, $false, $true, 1, $false, "47÷3=15, 2", 2) | Out-Null
$d.Content.Find.Execute("16÷3=5, 1", $true, $false, $false, $false, $false, $true, 1, $false, "48÷3=16, 0", 2) | Out-Null
$d.Content.Find.Execute("18÷4=4, 2", $true, $false, $false, $false, $false, $true, 1, $false, "47÷8=5, 7", 2) | Out-Null
$d.Content.Find.Execute("44÷5=8, 4", $true, $false, $false, $false, $false, $true, 1, $false, "51÷5=10, 1", 2) | Out-Null
$d.Content.Find.Execute("79÷2=39, 1", $true, $false, $false, $false, $false, $true, 1, $false, "85÷8=10, 5", 2) | Out-Null
$d.Content.Find.Execute("14÷8=1, 6", $true, $false, $false, $false, $false, $true, 1, $false, "10÷9=1, 1", 2) | Out-Null
$d.Content.Find.Execute("26÷7=3, 5", $true, $false, $false, $false, $false, $true, 1, $false, "19÷4=4, 3", 2) | Out-Null
$d.Content.Find.Execute("21÷6=3, 3", $true, $false, $false, $false, $false, $true, 1, $false, "64÷5=12, 4", 2) | Out-Null
